$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 3511  # was 3508
$ws.Range("F5").Value = 8223  # was 8219
$ws.Range("F7").Value = 86  # was 85
$ws.Range("F8").Value = 2172  # was 2171
$ws.Range("F10").Value = 185  # was 184
$ws.Range("F16").Value = 9  # was 8
$ws.Range("F19").Value = 1614  # was 1473
$ws.Range("F22").Value = 7193  # was 7188
$ws.Range("F24").Value = 55310  # was 55274
$ws.Range("F25").Value = 55310  # was 55275
$ws.Range("F26").Value = 4420  # was 4413
$ws.Range("F28").Value = 858  # was 856
$ws.Range("F33").Value = 2910  # was 2909
$ws.Range("F35").Value = 35  # was 33
$ws.Range("F38").Value = 1187  # was 1183
$ws.Range("F39").Value = 1117  # was 1106
$ws.Range("G43").Value = 68  # was 29.9
$ws.Range("F44").Value = 762  # was 761
$ws.Range("F45").Value = 155  # was 154
$ws.Range("F47").Value = 155  # was 154
$ws.Range("F49").Value = 38  # was 37
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 30  # was 29
$ws.Range("F10").Value = 50  # was 49
$ws.Range("F16").Value = 7464  # was 7463
$ws.Range("F17").Value = 104  # was 103
$ws.Range("F28").Value = 116  # was 114
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 2282  # was 2280
$ws.Range("F5").Value = 1540  # was 1538
$ws.Range("F8").Value = 2334  # was 2335
$ws.Range("F9").Value = 9331  # was 9329
$ws.Range("F10").Value = 1653  # was 1652
$ws.Range("F15").Value = 159  # was 155
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 8223  # was 8219
$ws.Range("F6").Value = 1653  # was 1652
$ws.Range("F9").Value = 86  # was 85
$ws.Range("F14").Value = 185  # was 184
$ws.Range("F17").Value = 1614  # was 1473
$ws.Range("F19").Value = 55310  # was 55275
$ws.Range("F22").Value = 858  # was 856
$ws.Range("F29").Value = 35  # was 33
$ws.Range("F31").Value = 1187  # was 1183
$ws.Range("F32").Value = 1117  # was 1106
$ws.Range("F33").Value = 104  # was 103
$ws.Range("G37").Value = 68  # was 29.9
$ws.Range("F38").Value = 762  # was 761
$ws.Range("F40").Value = 155  # was 154
$ws.Range("F42").Value = 116  # was 114
$ws.Range("F44").Value = 155  # was 154
$ws.Range("F46").Value = 38  # was 37
